# Update the dSF (column F) values for the steele_justin 2024 save-data sheet.
# This corresponds to a "repull data, push all data, mean calculation" update
# where the dSF column was recalculated for most rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new dSF (column F) value
$updates = @{
    2  = 4
    4  = 7
    5  = 4
    6  = 1
    8  = 4
    9  = -5
    10 = 1
    11 = -4
    12 = 8
    13 = 4
    14 = 0
    15 = 2
    17 = 2
    18 = -2
    21 = -1
    22 = -2
    23 = -5
    27 = -1
    28 = -2
    29 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
